$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume change (E) columns with refreshed data.
# D-column numeric-looking values are written as text (NumberFormat '@') so that
# trailing zeros / thousand-dot formatting (e.g. '0.0840', '2.20', '45.935.35')
# survive exactly as scraped, matching the original inline-string cell type.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.935.35'
$ws.Range("E2").Value = '  -0.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.599.66'
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.87'
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.89'
$ws.Range("E6").Value = '  -0.94%  '

$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.71'
$ws.Range("E10").Value = '  +0.85%  '

$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0840'
$ws.Range("E12").Value = '  +0.30%  '

$ws.Range("E13").Value = '  -3.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.002.68'
$ws.Range("E14").Value = '  +0.85%  '

$ws.Range("E15").Value = '  +1.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.605.68'
$ws.Range("E16").Value = '  +0.19%  '

$ws.Range("E17").Value = '  +0.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.77'
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '46.112.23'
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.74'
$ws.Range("E21").Value = '  +1.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.60'
$ws.Range("E22").Value = '  -2.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '289.69'
$ws.Range("E23").Value = '  +14.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.57'
$ws.Range("E24").Value = '  +1.98%  '

$ws.Range("E25").Value = '  +0.76%  '

$ws.Range("E26").Value = '  +2.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.49'
$ws.Range("E27").Value = '  +4.48%  '

$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.05'
$ws.Range("E29").Value = '  +0.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.71'
$ws.Range("E30").Value = '  +2.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.66'
$ws.Range("E31").Value = '  -2.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.20'
$ws.Range("E32").Value = '  -3.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.22'
$ws.Range("E33").Value = '  +2.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.62'
$ws.Range("E34").Value = '  -1.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.17'
$ws.Range("E35").Value = '  +4.10%  '

$ws.Range("E36").Value = '  -2.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0838'
$ws.Range("E37").Value = '  +1.72%  '

$ws.Range("E38").Value = '  -3.60%  '

$ws.Range("E39").Value = '  +4.03%  '

$ws.Range("E40").Value = '  +1.09%  '

$ws.Range("E41").Value = '  -2.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0328'
$ws.Range("E42").Value = '  +2.53%  '

$ws.Range("E43").Value = '  -2.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.36'
$ws.Range("E44").Value = '  +7.78%  '

$ws.Range("E45").Value = '  -4.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.117.66'
$ws.Range("E46").Value = '  +3.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '95.68'
$ws.Range("E47").Value = '  +5.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.998'
$ws.Range("E48").Value = '  -0.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.29'
$ws.Range("E49").Value = '  +0.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '109.10'
$ws.Range("E50").Value = '  +0.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.863.03'
$ws.Range("E51").Value = '  +0.86%  '
